$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 10
$ws.Range("E10").Value = 555
$ws.Range("F10").Value = 264
$ws.Range("H10").Value = 359

# Row 11
$ws.Range("E11").Value = 355
$ws.Range("F11").Value = 190
$ws.Range("H11").Value = 255

# Row 14
$ws.Range("E14").Value = 124
$ws.Range("F14").Value = 65
$ws.Range("H14").Value = 99

# Row 15
$ws.Range("E15").Value = 167
$ws.Range("F15").Value = 72
$ws.Range("H15").Value = 122

# Row 23
$ws.Range("E23").Value = 200
$ws.Range("F23").Value = 94
$ws.Range("H23").Value = 145

# Row 28
$ws.Range("E28").Value = 198

# Row 35
$ws.Range("F35").Value = 93
$ws.Range("H35").Value = 120

# Row 41
$ws.Range("E41").Value = 388
